# Refresh generated Liga Classica datasets:
#  - "Geral": update each team's total score (column B) with the
#    freshly recomputed values from the automated notebook run.
#  - "Mes - Janeiro": this tab is a ranking of the same teams sorted
#    descending by score, so it is rebuilt from the refreshed totals.

$wb = $excel.ActiveWorkbook

# ---- 1) Update the "Geral" totals (column B), in sheet row order --------
$geral = $wb.Worksheets.Item("Geral")

$geralRows = 2, 4, 5, 6, 7, 8, 10, 11, 12, 13, 14, 15, 17, 18, 20, 21
$geralNewScores = 56.69, 42.06, 62.07, 83.2, 49.36, 67.66, 57.26, 54.05, 43.56, 43.96, 58.26, 55.66, 18.49, 56.86, 45.59, 63.76

for ($i = 0; $i -lt $geralRows.Length; $i++) {
    $geral.Cells.Item($geralRows[$i], 2).Value = $geralNewScores[$i]
}

# ---- 2) Rebuild the "Mes - Janeiro" ranking ------------------------------
# Ranking order: teams sorted by their (refreshed) Geral score, descending.
$jan = $wb.Worksheets.Item("Mes - Janeiro")

$ranking = "Esquadrão Gazembrino", "GaúchoDaFronteira F.C", "SC 100 Sono", "Texas Club 2026", "bugredasmissões", "Doug Leal F.C", "lsauer fc", "GrioTeam", "Pontaç0 F.C.", "Arran Katoko FC", "Medonho´s F.C.", "Grêmio_Campeão_LA_27", "FBC Colorado", "SC ÉoINTER!", "GE Bebum", "La Primeira Patada Es Nuestra", "JV5 Tricolor Gaúcho", "C R Juvenal", "Pepe Leal FC", "NHU PORÃ SAF."
$rankingScores = 83.2, 67.66, 67.16, 63.76, 62.76, 62.07, 58.26, 57.26, 56.86, 56.69, 55.66, 54.05, 49.36, 45.59, 44.65, 43.96, 43.56, 42.06, 18.49, 0

for ($i = 0; $i -lt $ranking.Length; $i++) {
    $row = $i + 2
    $jan.Cells.Item($row, 1).Value = $ranking[$i]
    $jan.Cells.Item($row, 2).Value = $rankingScores[$i]
}
